# Auto-generated Excel COM-interop script to apply Lamia_Profits.xlsx numeric updates
# across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets (columns H-N, various rows).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 137
$ws.Range("H137").Value = 29414852
$ws.Range("I137").Value = 55557684
$ws.Range("J137").Value = 4164.75
$ws.Range("K137").Value = 166673052
$ws.Range("L137").Value = 12494.25
$ws.Range("M137").Value = -166670502
$ws.Range("N137").Value = -17594.25

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 3226.2222
$ws.Range("I45").Value = 2158.4666
$ws.Range("J45").Value = 8565
$ws.Range("K45").Value = 2158.4666
$ws.Range("L45").Value = 8565
$ws.Range("M45").Value = -1781.4666
$ws.Range("N45").Value = -9319

# Row 61
$ws.Range("H61").Value = 3954.1475
$ws.Range("I61").Value = 3885
$ws.Range("J61").Value = 5994
$ws.Range("K61").Value = 3885
$ws.Range("L61").Value = 5994
$ws.Range("M61").Value = -3673
$ws.Range("N61").Value = -6418

# Row 88
$ws.Range("H88").Value = 3114.3157
$ws.Range("I88").Value = 3071
$ws.Range("J88").Value = 3139.5833
$ws.Range("K88").Value = 3071
$ws.Range("L88").Value = 3139.5833
$ws.Range("M88").Value = -2665
$ws.Range("N88").Value = -3951.5833

# Row 91
$ws.Range("H91").Value = 3114.3157
$ws.Range("I91").Value = 3071
$ws.Range("J91").Value = 3139.5833
$ws.Range("K91").Value = 3071
$ws.Range("L91").Value = 3139.5833
$ws.Range("M91").Value = -1667
$ws.Range("N91").Value = -5947.5833

# Row 136
$ws.Range("H136").Value = 3954.1475
$ws.Range("I136").Value = 3885
$ws.Range("J136").Value = 5994
$ws.Range("K136").Value = 11655
$ws.Range("L136").Value = 17982
$ws.Range("M136").Value = -9105
$ws.Range("N136").Value = -23082

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2643.7
$ws.Range("I86").Value = 2611.8845
$ws.Range("J86").Value = 2850.5
$ws.Range("K86").Value = 2611.8845
$ws.Range("L86").Value = 2850.5
$ws.Range("M86").Value = -1488.8845
$ws.Range("N86").Value = -5096.5

# Row 89
$ws.Range("H89").Value = 2643.7
$ws.Range("I89").Value = 2611.8845
$ws.Range("J89").Value = 2850.5
$ws.Range("K89").Value = 13059.4225
$ws.Range("L89").Value = 14252.5
$ws.Range("M89").Value = -7443.422500000001
$ws.Range("N89").Value = -25484.5

# Row 94
$ws.Range("H94").Value = 1400.0667
$ws.Range("I94").Value = 1375.0834
$ws.Range("J94").Value = 1500
$ws.Range("K94").Value = 1375.0834
$ws.Range("L94").Value = 1500
$ws.Range("M94").Value = -924.0834
$ws.Range("N94").Value = -2402

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 4145.36
$ws.Range("I58").Value = 1939.4445
$ws.Range("J58").Value = 9817.714
$ws.Range("K58").Value = 1939.4445
$ws.Range("L58").Value = 9817.714
$ws.Range("M58").Value = -1736.4445
$ws.Range("N58").Value = -10223.714

# Row 99
$ws.Range("H99").Value = 3498.25
$ws.Range("I99").Value = 3498
$ws.Range("J99").Value = 3499
$ws.Range("K99").Value = 3498
$ws.Range("L99").Value = 3499
$ws.Range("M99").Value = -2000
$ws.Range("N99").Value = -6495

# Row 126
$ws.Range("H126").Value = 3498.25
$ws.Range("I126").Value = 3498
$ws.Range("J126").Value = 3499
$ws.Range("K126").Value = 10494
$ws.Range("L126").Value = 3499
$ws.Range("M126").Value = -8024
$ws.Range("N126").Value = -15437

# Row 136
$ws.Range("H136").Value = 4145.36
$ws.Range("I136").Value = 1939.4445
$ws.Range("J136").Value = 9817.714
$ws.Range("K136").Value = 5818.333500000001
$ws.Range("L136").Value = 29453.142
$ws.Range("M136").Value = -3268.333500000001
$ws.Range("N136").Value = -34553.142

$ws = $wb.Worksheets.Item("CUL")
# Row 38
$ws.Range("H38").Value = 36.1875
$ws.Range("I38").Value = 57.666668
$ws.Range("J38").Value = 23.3
$ws.Range("K38").Value = 173.000004
$ws.Range("L38").Value = 69.90000000000001
$ws.Range("M38").Value = 173.999996
$ws.Range("N38").Value = -763.9

# Row 42
$ws.Range("H42").Value = 7500
$ws.Range("I42").Value = 5000
$ws.Range("J42").Value = 10000
$ws.Range("K42").Value = 15000
$ws.Range("L42").Value = 30000
$ws.Range("M42").Value = -14466
$ws.Range("N42").Value = -31068

# Row 92
$ws.Range("H92").Value = 1226.7
$ws.Range("I92").Value = 435.375
$ws.Range("J92").Value = 4392
$ws.Range("K92").Value = 1306.125
$ws.Range("L92").Value = 13176
$ws.Range("M92").Value = -58.125
$ws.Range("N92").Value = -15672

# Row 114
$ws.Range("H114").Value = 1349.8
$ws.Range("I114").Value = 795
$ws.Range("J114").Value = 1719.6666
$ws.Range("K114").Value = 2385
$ws.Range("L114").Value = 5158.9998
$ws.Range("M114").Value = 869
$ws.Range("N114").Value = -11666.9998

# Row 132
$ws.Range("H132").Value = 4234.227
$ws.Range("I132").Value = 2805.0908
$ws.Range("J132").Value = 5663.364
$ws.Range("K132").Value = 25245.8172
$ws.Range("L132").Value = 50970.276
$ws.Range("M132").Value = -22715.8172
$ws.Range("N132").Value = -56030.276

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 560.26666
$ws.Range("I2").Value = 95
$ws.Range("J2").Value = 1092
$ws.Range("K2").Value = 95
$ws.Range("L2").Value = 1092
$ws.Range("M2").Value = 18
$ws.Range("N2").Value = -1318

# Row 40
$ws.Range("H40").Value = 8000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 8000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 8000
$ws.Range("N40").Value = -8302

# Row 80
$ws.Range("H80").Value = 10150.4
$ws.Range("I80").Value = 8100
$ws.Range("J80").Value = 12200.8
$ws.Range("K80").Value = 8100
$ws.Range("L80").Value = 12200.8
$ws.Range("M80").Value = -7102
$ws.Range("N80").Value = -14196.8

# Row 83
$ws.Range("H83").Value = 10150.4
$ws.Range("I83").Value = 8100
$ws.Range("J83").Value = 12200.8
$ws.Range("K83").Value = 40500
$ws.Range("L83").Value = 61004
$ws.Range("M83").Value = -35508
$ws.Range("N83").Value = -70988

# Row 113
$ws.Range("H113").Value = 1247.5
$ws.Range("I113").Value = 1247.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1247.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 922.5
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 3512.875
$ws.Range("I16").Value = 2585.8572
$ws.Range("J16").Value = 10002
$ws.Range("K16").Value = 2585.8572
$ws.Range("L16").Value = 10002
$ws.Range("M16").Value = -2415.8572
$ws.Range("N16").Value = -10342

# Row 22
$ws.Range("H22").Value = 9274
$ws.Range("I22").Value = 5483
$ws.Range("J22").Value = 11341.818
$ws.Range("K22").Value = 5483
$ws.Range("L22").Value = 11341.818
$ws.Range("M22").Value = -5188
$ws.Range("N22").Value = -11931.818

# Row 27
$ws.Range("H27").Value = 9274
$ws.Range("I27").Value = 5483
$ws.Range("J27").Value = 11341.818
$ws.Range("K27").Value = 5483
$ws.Range("L27").Value = 11341.818
$ws.Range("M27").Value = -5376
$ws.Range("N27").Value = -11555.818

# Row 100
$ws.Range("H100").Value = 10626.381
$ws.Range("I100").Value = 7684.048
$ws.Range("J100").Value = 13568.714
$ws.Range("K100").Value = 7684.048
$ws.Range("L100").Value = 13568.714
$ws.Range("M100").Value = -7143.048
$ws.Range("N100").Value = -14650.714

# Row 122
$ws.Range("H122").Value = 5499.9165
$ws.Range("I122").Value = 3512
$ws.Range("J122").Value = 9475.75
$ws.Range("K122").Value = 10536
$ws.Range("L122").Value = 28427.25
$ws.Range("M122").Value = -8086
$ws.Range("N122").Value = -33327.25

# Row 132
$ws.Range("H132").Value = 10385.308
$ws.Range("I132").Value = 2624.75
$ws.Range("J132").Value = 13834.444
$ws.Range("K132").Value = 7874.25
$ws.Range("L132").Value = 41503.33199999999
$ws.Range("M132").Value = -5344.25
$ws.Range("N132").Value = -46563.33199999999

$ws = $wb.Worksheets.Item("WVR")
# Row 5
$ws.Range("H5").Value = 12513750
$ws.Range("I5").Value = 17499
$ws.Range("J5").Value = 25010000
$ws.Range("K5").Value = 17499
$ws.Range("L5").Value = 25010000
$ws.Range("M5").Value = -17387
$ws.Range("N5").Value = -25010224

# Row 41
$ws.Range("H41").Value = 13827.714
$ws.Range("I41").Value = 10500
$ws.Range("J41").Value = 15158.8
$ws.Range("K41").Value = 10500
$ws.Range("L41").Value = 15158.8
$ws.Range("M41").Value = -10110
$ws.Range("N41").Value = -15938.8

# Row 100
$ws.Range("H100").Value = 963.5263
$ws.Range("I100").Value = 829.5
$ws.Range("J100").Value = 1193.2858
$ws.Range("K100").Value = 1659
$ws.Range("L100").Value = 2386.5716
$ws.Range("M100").Value = -1118
$ws.Range("N100").Value = -3468.5716

# Row 106
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
